# "autoIt and trace log" - update the Google-sheet trace log with the
# latest automated test run's results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Google")

# Test case #1 (row 3) no longer produced a first-result / URL for this run.
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""

# Test case #4 (row 6) flipped from pass to fail.
$ws.Range("D6").Value = "fail"
